$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-04 18:25:04"

# Refresh the "取得日時" (fetched at) timestamp for every data row (2-15)
for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}

# Row 13 and Row 14 swapped places in the source feed: carry over their
# title (B) and URL (F) values so each row now shows the other's listing.
$b13 = $ws.Range("B13").Value()
$f13 = $ws.Range("F13").Value()
$b14 = $ws.Range("B14").Value()
$f14 = $ws.Range("F14").Value()

$ws.Range("B13").Value = $b14
$ws.Range("F13").Value = $f14
$ws.Range("B14").Value = $b13
$ws.Range("F14").Value = $f13
